$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "中文在线"
$ws.Range("B2").Value = "协鑫集成"
$ws.Range("C2").Value = "协鑫集成"
$ws.Range("A3").Value = "光线传媒"
$ws.Range("B3").Value = "中文在线"
$ws.Range("C3").Value = "万向钱潮"
$ws.Range("A4").Value = "捷成股份"
$ws.Range("B4").Value = "光线传媒"
$ws.Range("C4").Value = "巨力索具"
$ws.Range("A5").Value = "协鑫集成"
$ws.Range("B5").Value = "捷成股份"
$ws.Range("C5").Value = "光线传媒"
$ws.Range("A6").Value = "巨力索具"
$ws.Range("B6").Value = "巨力索具"
$ws.Range("C6").Value = "博纳影业"
$ws.Range("A7").Value = "掌阅科技"
$ws.Range("B7").Value = "博纳影业"
$ws.Range("C7").Value = "中文在线"
$ws.Range("A8").Value = "博纳影业"
$ws.Range("B8").Value = "掌阅科技"
$ws.Range("C8").Value = "横店影视"
$ws.Range("A9").Value = "横店影视"
$ws.Range("B9").Value = "万向钱潮"
$ws.Range("C9").Value = "嘉美包装"
$ws.Range("A10").Value = "万向钱潮"
$ws.Range("B10").Value = "杉杉股份"
$ws.Range("C10").Value = "捷成股份"
$ws.Range("A11").Value = "蓝色光标"
$ws.Range("B11").Value = "湖南白银"
$ws.Range("C11").Value = "中超控股"
$ws.Range("A12").Value = "视觉中国"
$ws.Range("B12").Value = "利欧股份"
$ws.Range("C12").Value = "利欧股份"
$ws.Range("A13").Value = "利欧股份"
$ws.Range("B13").Value = "视觉中国"
$ws.Range("C13").Value = "浙文互联"
$ws.Range("A14").Value = "大位科技"
$ws.Range("B14").Value = "蓝色光标"
$ws.Range("C14").Value = "视觉中国"
$ws.Range("A15").Value = "长飞光纤"
$ws.Range("B15").Value = "五洲新春"
$ws.Range("C15").Value = "掌阅科技"
$ws.Range("A16").Value = "浙文互联"
$ws.Range("B16").Value = "横店影视"
$ws.Range("C16").Value = "大位科技"
$ws.Range("A17").Value = "五洲新春"
$ws.Range("B17").Value = "数据港"
$ws.Range("C17").Value = "航天发展"
$ws.Range("A18").Value = "数据港"
$ws.Range("B18").Value = "中国电影"
$ws.Range("C18").Value = "五洲新春"
$ws.Range("A19").Value = "中超控股"
$ws.Range("B19").Value = "大位科技"
$ws.Range("C19").Value = "蓝色光标"
$ws.Range("A20").Value = "嘉美包装"
$ws.Range("B20").Value = "浙文互联"
$ws.Range("C20").Value = "百川股份"
$ws.Range("A21").Value = "杉杉股份"
$ws.Range("B21").Value = "特发信息"
$ws.Range("C21").Value = "长飞光纤"

Write-Output "Applied HotStock_Top20 cell updates"
